# Update the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# timestamps for the a7ed6a24... entry (row 4) on both the zh-cn and de-de
# report sheets, per the latest handback run.
#
# zh-cn sheet:
#   Correspond Handoff Datetime : 2016-03-12 06:13:46 -> 2016-03-12 06:14:30
#   Correspond Handback DateTime: 2016-03-12 06:14:06 -> 2016-03-12 06:14:47
#
# de-de sheet:
#   Correspond Handoff Datetime : 2016-03-12 06:13:49 -> 2016-03-12 06:14:33
#   Correspond Handback DateTime: 2016-03-12 06:14:11 -> 2016-03-12 06:14:53

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4").Value = "2016-03-12 06:14:30"
$wsZhCn.Range("E5").Value = "2016-03-12 06:14:30"
$wsZhCn.Range("H4").Value = "2016-03-12 06:14:47"
$wsZhCn.Range("H5").Value = "2016-03-12 06:14:47"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4").Value = "2016-03-12 06:14:33"
$wsDeDe.Range("E5").Value = "2016-03-12 06:14:33"
$wsDeDe.Range("H4").Value = "2016-03-12 06:14:53"
$wsDeDe.Range("H5").Value = "2016-03-12 06:14:53"
